$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 182.07692
$ws.Range("I33").Value = 163.83333
$ws.Range("K33").Value = 163.83333
$ws.Range("M33").Value = 65.16667000000001
$ws.Range("H63").Value = 26999.5
$ws.Range("J63").Value = 26999.5
$ws.Range("L63").Value = 26999.5
$ws.Range("N63").Value = -28247.5
$ws.Range("H66").Value = 26999.5
$ws.Range("J66").Value = 26999.5
$ws.Range("L66").Value = 80998.5
$ws.Range("N66").Value = -87238.5
$ws.Range("H98").Value = 5819.1636
$ws.Range("I98").Value = 4551.6553
$ws.Range("J98").Value = 7232.923
$ws.Range("K98").Value = 4551.6553
$ws.Range("L98").Value = 7232.923
$ws.Range("M98").Value = -3053.6553
$ws.Range("N98").Value = -10228.923
$ws.Range("H106").Value = 1313
$ws.Range("I106").Value = 916.25
$ws.Range("J106").Value = 2900
$ws.Range("K106").Value = 916.25
$ws.Range("L106").Value = 2900
$ws.Range("M106").Value = -285.25
$ws.Range("N106").Value = -4162
$ws.Range("H112").Value = 1260.1562
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 1297.541
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 3892.623
$ws.Range("M112").Value = -392
$ws.Range("N112").Value = -6108.623
$ws.Range("H121").Value = 1000
$ws.Range("I121").Value = 100
$ws.Range("J121").Value = 1300
$ws.Range("K121").Value = 300
$ws.Range("L121").Value = 3900
$ws.Range("M121").Value = 1447
$ws.Range("N121").Value = -7394
$ws.Range("H122").Value = 5819.1636
$ws.Range("I122").Value = 4551.6553
$ws.Range("J122").Value = 7232.923
$ws.Range("K122").Value = 13654.9659
$ws.Range("L122").Value = 21698.769
$ws.Range("M122").Value = -11204.9659
$ws.Range("N122").Value = -26598.769
$ws.Range("H129").Value = 842.45
$ws.Range("J129").Value = 861.90625
$ws.Range("L129").Value = 2585.71875
$ws.Range("N129").Value = -12585.71875
$ws.Range("H131").Value = 2910.625
$ws.Range("I131").Value = 2169.5454
$ws.Range("J131").Value = 4541
$ws.Range("K131").Value = 6508.6362
$ws.Range("L131").Value = 13623
$ws.Range("M131").Value = -1468.6362
$ws.Range("N131").Value = -23703
$ws.Range("H138").Value = 2979.7827
$ws.Range("I138").Value = 2042.5
$ws.Range("J138").Value = 3037.4614
$ws.Range("K138").Value = 6127.5
$ws.Range("L138").Value = 9112.3842
$ws.Range("M138").Value = -987.5
$ws.Range("N138").Value = -19392.3842

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 36000
$ws.Range("J103").Value = 36000
$ws.Range("L103").Value = 36000
$ws.Range("N103").Value = -38344
$ws.Range("H107").Value = 852.40625
$ws.Range("I107").Value = 690.3043
$ws.Range("J107").Value = 1266.6666
$ws.Range("K107").Value = 690.3043
$ws.Range("L107").Value = 1266.6666
$ws.Range("M107").Value = 1229.6957
$ws.Range("N107").Value = -5106.6666
$ws.Range("H112").Value = 29990
$ws.Range("J112").Value = 29990
$ws.Range("L112").Value = 29990
$ws.Range("N112").Value = -32944
$ws.Range("H125").Value = 41830
$ws.Range("J125").Value = 41830
$ws.Range("L125").Value = 41830
$ws.Range("N125").Value = -51670
$ws.Range("H130").Value = 41920
$ws.Range("J130").Value = 41920
$ws.Range("L130").Value = 41920
$ws.Range("N130").Value = -51960
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H135").Value = 47249.855
$ws.Range("J135").Value = 47249.855
$ws.Range("L135").Value = 47249.855
$ws.Range("N135").Value = -57389.855
$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200
$ws.Range("H138").Value = 40956.523
$ws.Range("J138").Value = 40956.523
$ws.Range("L138").Value = 40956.523
$ws.Range("N138").Value = -51236.523
$ws.Range("H140").Value = 50545.043
$ws.Range("J140").Value = 50545.043
$ws.Range("L140").Value = 50545.043
$ws.Range("N140").Value = -60905.043

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8384.677
$ws.Range("I31").Value = 4077
$ws.Range("J31").Value = 11400.05
$ws.Range("K31").Value = 4077
$ws.Range("L31").Value = 11400.05
$ws.Range("M31").Value = -3782
$ws.Range("N31").Value = -11990.05
$ws.Range("H34").Value = 8384.677
$ws.Range("I34").Value = 4077
$ws.Range("J34").Value = 11400.05
$ws.Range("K34").Value = 4077
$ws.Range("L34").Value = 11400.05
$ws.Range("M34").Value = -3875
$ws.Range("N34").Value = -11804.05
$ws.Range("H94").Value = 1796.3125
$ws.Range("J94").Value = 1936.8334
$ws.Range("L94").Value = 1936.8334
$ws.Range("N94").Value = -2838.8334
$ws.Range("H134").Value = 5721.2593
$ws.Range("I134").Value = 6945.8887
$ws.Range("J134").Value = 3272
$ws.Range("K134").Value = 20837.6661
$ws.Range("L134").Value = 9816
$ws.Range("M134").Value = -18302.6661
$ws.Range("N134").Value = -14886

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 25050566
$ws.Range("I107").Value = 469.85715
$ws.Range("J107").Value = 38539080
$ws.Range("K107").Value = 1409.57145
$ws.Range("L107").Value = 115617240
$ws.Range("M107").Value = 510.4285500000001
$ws.Range("N107").Value = -115621080
$ws.Range("H131").Value = 777.86
$ws.Range("J131").Value = 824.5824
$ws.Range("L131").Value = 2473.7472
$ws.Range("N131").Value = -12553.7472

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 11491.167
$ws.Range("I43").Value = 1426.3334
$ws.Range("J43").Value = 28265.889
$ws.Range("K43").Value = 1426.3334
$ws.Range("L43").Value = 28265.889
$ws.Range("M43").Value = -1275.3334
$ws.Range("N43").Value = -28567.889
$ws.Range("H46").Value = 27816.455
$ws.Range("J46").Value = 28073.1
$ws.Range("L46").Value = 28073.1
$ws.Range("N46").Value = -28385.1
$ws.Range("H57").Value = 34860.43
$ws.Range("J57").Value = 34860.43
$ws.Range("L57").Value = 34860.43
$ws.Range("N57").Value = -36500.43
$ws.Range("H80").Value = 10873288
$ws.Range("I80").Value = 20836810
$ws.Range("J80").Value = 3991.4546
$ws.Range("K80").Value = 20836810
$ws.Range("L80").Value = 3991.4546
$ws.Range("M80").Value = -20835812
$ws.Range("N80").Value = -5987.4546
$ws.Range("H83").Value = 10873288
$ws.Range("I83").Value = 20836810
$ws.Range("J83").Value = 3991.4546
$ws.Range("K83").Value = 104184050
$ws.Range("L83").Value = 19957.273
$ws.Range("M83").Value = -104179058
$ws.Range("N83").Value = -29941.273
$ws.Range("H97").Value = 1054.55
$ws.Range("I97").Value = 1053.3334
$ws.Range("J97").Value = 1058.2
$ws.Range("K97").Value = 1053.3334
$ws.Range("L97").Value = 1058.2
$ws.Range("M97").Value = -557.3334
$ws.Range("N97").Value = -2050.2
$ws.Range("H132").Value = 3517.625
$ws.Range("I132").Value = 1662.1818
$ws.Range("J132").Value = 7599.6
$ws.Range("K132").Value = 4986.5454
$ws.Range("L132").Value = 22798.8
$ws.Range("M132").Value = -2456.5454
$ws.Range("N132").Value = -27858.8

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 993
$ws.Range("I16").Value = 992.3333
$ws.Range("J16").Value = 999
$ws.Range("K16").Value = 992.3333
$ws.Range("L16").Value = 999
$ws.Range("M16").Value = -822.3333
$ws.Range("N16").Value = -1339
$ws.Range("H46").Value = 1407.0741
$ws.Range("I46").Value = 833.26666
$ws.Range("J46").Value = 2124.3333
$ws.Range("K46").Value = 833.26666
$ws.Range("L46").Value = 2124.3333
$ws.Range("M46").Value = -645.26666
$ws.Range("N46").Value = -2500.3333

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 16405.818
$ws.Range("J54").Value = 16405.818
$ws.Range("L54").Value = 16405.818
$ws.Range("N54").Value = -17445.818
$ws.Range("H64").Value = 25977.777
$ws.Range("J64").Value = 25977.777
$ws.Range("L64").Value = 25977.777
$ws.Range("N64").Value = -26473.777
$ws.Range("H67").Value = 25977.777
$ws.Range("J67").Value = 25977.777
$ws.Range("L67").Value = 25977.777
$ws.Range("N67").Value = -27693.777
$ws.Range("H96").Value = 38866564
$ws.Range("I96").Value = 48119930
$ws.Range("J96").Value = 2429.8
$ws.Range("K96").Value = 48119930
$ws.Range("L96").Value = 2429.8
$ws.Range("M96").Value = -48118557
$ws.Range("N96").Value = -5175.8
$ws.Range("H122").Value = 3822.8918
$ws.Range("I122").Value = 2637.35
$ws.Range("J122").Value = 5217.647
$ws.Range("K122").Value = 7912.049999999999
$ws.Range("L122").Value = 15652.941
$ws.Range("M122").Value = -5462.049999999999
$ws.Range("N122").Value = -20552.941
